$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Update the "last saved / displayed" date fields (4/11/2019 -> 4/12/2019) ---
# These datetimeFigureOut fields live on the notes master, the slide master and
# every slide layout (not on the slide itself, since this deck has no notes
# pages and slide1 does not carry its own date placeholder).
function Set-DateFieldText {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "4/11/2019") {
                    $tr.Text = "4/12/2019"
                }
            }
        }
    }
}

Set-DateFieldText $p.NotesMaster.Shapes
Set-DateFieldText $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Set-DateFieldText $p.SlideMaster.CustomLayouts.Item($li).Shapes
}

# --- Add the new dashed connector shape to slide 1 ---
# The original deck had shapes deleted earlier in its life, so PowerPoint's
# internal id counter is well ahead of the visible shape count; the shape
# that was added in the real edit picked up id=52. Burn through the
# intervening ids with throwaway shapes (deleting each one) until the
# *next* shape created lands exactly on id 52 - that one becomes our real
# shape instead of being thrown away.
$dummies = @()
$shp = $null
while ($true) {
    $probe = $s.Shapes.AddLine(0, 0, 1, 1)
    if ($probe.Id -ge 52) {
        $shp = $probe
        break
    }
    $dummies += $probe
}

$shp.Left = 781.7424409448819
$shp.Top = 438.81464566929134
$shp.Width = 0
$shp.Height = 32.22629921259843
$shp.Name = "Straight Connector 51"

$shp.Line.ForeColor.RGB = 0xC07000
$shp.Line.Weight = 1.5
$shp.Line.DashStyle = 9

$shp.Shadow.Visible = $false

foreach ($d in $dummies) {
    $d.Delete()
}
